# Progress-Report.xlsx update:
#  - "Uploaded most recent test cases & updated progress report"
#  - Fill in PERSON-IN-CHARGE / REMARKS on "User Interface"
#  - Fill in % COMPLETE / DONE / PERSON-IN-CHARGE on "Back end"
#  - Fill in % COMPLETE / DONE / PERSON-IN-CHARGE on "Test Cases"
#  - Refresh view/selection state and normalize row heights

$wb = $excel.ActiveWorkbook

$wsUI   = $wb.Worksheets.Item(1)   # User Interface
$wsBack = $wb.Worksheets.Item(2)   # Back end
$wsTest = $wb.Worksheets.Item(3)   # Test Cases
$wsPaper = $wb.Worksheets.Item(4)  # Paperworks

# ---------------------------------------------------------------------------
# User Interface sheet (Table14: B=PAGE, C=%COMPLETE, D=DONE, E=REMARKS, F=PERSON-IN-CHARGE)
# New remarks/names must be entered in this exact order so the shared-string
# table picks up new unique strings in the same sequence as the source edit.
# ---------------------------------------------------------------------------
$wsUI.Range("E7").Value = "No functionality"
$wsUI.Range("E9").Value = "Ask Miguel / Daniel?, Fix layout, match design with other pages"
$wsUI.Range("F4").Value = "Daniel Garcia"
$wsUI.Range("F6").Value = "Kenywil Tiu"
$wsUI.Range("F3").Value = "Margaret Avancena"
$wsUI.Range("F7").Value = "Margaret Avancena"
$wsUI.Range("F8").Value = "Rheygine Medel"
$wsUI.Range("F9").Value = "Daniel Garcia"

# Remove the stale remark that used to sit on row 6
$wsUI.Range("E6").ClearContents()

# Progress updates
$wsUI.Range("C6").Value = 1
$wsUI.Range("C7").Value = 1
$wsUI.Range("C8").Value = 1

# ---------------------------------------------------------------------------
# Back end sheet (Table1: B=FUNCTIONS, C=USER STORY#, D=ITERATION,
#                 E=%COMPLETE, F=DONE, G=REMARKS, H=PERSON-IN-CHARGE)
# ---------------------------------------------------------------------------
$wsBack.Range("E12").Value = 0

$wsBack.Range("E14").Value = 1
$wsBack.Range("H14").Value = "Miguel Manalac"

$wsBack.Range("E15").Value = 1
$wsBack.Range("H15").Value = "Joseph Ongsingco"

$wsBack.Range("E16").Value = 1
$wsBack.Range("H16").Value = "Carlo Bautista"

# ---------------------------------------------------------------------------
# Test Cases sheet (Table13: B=USER STORY, C=ITERATION, D=%COMPLETE,
#                   E=DONE, F=REMARKS, G=PERSON-IN-CHARGE)
# ---------------------------------------------------------------------------
$wsTest.Range("D14").Value = 1
$wsTest.Range("D14").NumberFormat = "0%"
$wsTest.Range("G14").Value = "Brandon Partosa"

$wsTest.Range("D15").Value = 1
$wsTest.Range("D15").NumberFormat = "0%"
$wsTest.Range("G15").Value = "Ralph Chua"

$wsTest.Range("D16").Value = 1
$wsTest.Range("D16").NumberFormat = "0%"
$wsTest.Range("G16").Value = "Edgar Dimanarig"

# ---------------------------------------------------------------------------
# Normalize row heights (re-saved workbook rounds every row height to the
# nearest whole point).
# ---------------------------------------------------------------------------
$wsUI.Rows.Item("1:29").RowHeight = 22
$wsTest.Rows.Item("1:21").RowHeight = 22
$wsPaper.Rows.Item("1:13").RowHeight = 22

$wsBack.Rows.Item("1:5").RowHeight = 22
$wsBack.Rows.Item(6).RowHeight = 36
$wsBack.Rows.Item(7).RowHeight = 22
$wsBack.Rows.Item(8).RowHeight = 34
$wsBack.Rows.Item(9).RowHeight = 43
$wsBack.Rows.Item(10).RowHeight = 37
$wsBack.Rows.Item(11).RowHeight = 40
$wsBack.Rows.Item(12).RowHeight = 40
$wsBack.Rows.Item(13).RowHeight = 35
$wsBack.Rows.Item(14).RowHeight = 38
$wsBack.Rows.Item(15).RowHeight = 33
$wsBack.Rows.Item(16).RowHeight = 50
$wsBack.Rows.Item(17).RowHeight = 31
$wsBack.Rows.Item(18).RowHeight = 33
$wsBack.Rows.Item(19).RowHeight = 48
$wsBack.Rows.Item(20).RowHeight = 33
$wsBack.Rows.Item(21).RowHeight = 33

# ---------------------------------------------------------------------------
# View / selection state: move the active selection on each sheet, then
# finish with Paperworks as the active tab (matches the saved workbook view).
# ---------------------------------------------------------------------------
$wsUI.Activate()
$wsUI.Range("E22").Select()

$wsBack.Activate()
$wsBack.Range("C14:H14").Select()

$wsTest.Activate()
$wsTest.Range("G15").Select()

$wsPaper.Activate()
$wsPaper.Range("C9").Select()
